# Vwap.Calc.xlsx: rename the "index" column to "i" and switch it from a
# 1-based to a 0-based counter on the "VWAP" sheet (table `testdata7`).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VWAP")

# Header: "index" -> "i" (updates the table column name + shared string).
$ws.Range("A1").Value = "i"

# Data rows: values were 1..391, now 0..390 (shift down by one).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Column A narrows to fit the shorter header/values.
$ws.Columns.Item(1).ColumnWidth = 3.17
